# Generate Report for Handback
# Row 7 ("05b7e7f0-c301-4859-bace-a2dbe4497cf1.md") on both the zh-cn and
# de-de worksheets now has a handback result: the handback file showed up,
# but it is not the latest version of the source doc, so the report fills
# in the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" / "Error Detail" columns (I, J, K, P) for that row.

$wb = $excel.ActiveWorkbook

$latestSourceUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/0b8d87ade864b8f65025e64cdd79f4d95f82dc80/e2e/05b7e7f0-c301-4859-bace-a2dbe4497cf1.md"
$targetFileDisplay = "05b7e7f0-c301-4859-bace-a2dbe4497cf1.md"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/97ea3fe74ffc28e226762d1784e2876994a3c901/e2e/05b7e7f0-c301-4859-bace-a2dbe4497cf1.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/0b8d87ade864b8f65025e64cdd79f4d95f82dc80/e2e/05b7e7f0-c301-4859-bace-a2dbe4497cf1.md."

# zh-cn worksheet, row 7
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("J7").Value = "05b7e7f0-c301-4859-bace-a2dbe4497cf1.0569686849b2f666ed787bf98b990950057d7dc2.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-14 17:13:35"
$wsZh.Range("P7").Value = $errorDetail
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $latestSourceUrl, "", "", $targetFileDisplay)

# de-de worksheet, row 7
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("J7").Value = "05b7e7f0-c301-4859-bace-a2dbe4497cf1.0569686849b2f666ed787bf98b990950057d7dc2.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-14 17:13:46"
$wsDe.Range("P7").Value = $errorDetail
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $latestSourceUrl, "", "", $targetFileDisplay)
